$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = 16.485934880815449
$ws.Range("C2").Value = 11.310938574386626
$ws.Range("D2").Value = 13.177117714156338
$ws.Range("E2").Value = -0.17549645616645648

# Row 3 updates (C3 cleared, D3 newly populated)
$ws.Range("B3").Value = 34.039170959114387
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 2.287950224122568
$ws.Range("E3").Value = 2.7539060664816475

# Update the selection to match the new range used
$ws.Range("B1:E3").Select()
